$p = $ppt.ActivePresentation

$tableSlides = @(14, 15, 16)
foreach ($slideIndex in $tableSlides) {
    $s = $p.Slides.Item($slideIndex)
    $sh = $s.Shapes.Item(1)
    $sh.Table.ApplyStyle("{A466CCF1-3076-4506-B1A9-68C3800A951D}")
}
